# Generate Report for Handoff
# Re-generates the localization status report for the row corresponding to
# "6936653b-c66c-4b4e-9a9d-7e7c21f24e07.md", updating the handoff/generate
# timestamps and marking the handoff Priority/type as "ht" on the zh-cn and
# de-de worksheets.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 11, 13, 14)

# --- "Overview" sheet: column G = "Latest HO Xliff Generate Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-30 20:23:34"
}

# --- "zh-cn" sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-30 20:23:29"
}

# --- "de-de" sheet: column E = "Priority", column H = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-30 20:23:34"
}
